$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for "line7" and "line8" right after the "line6" row (row 8),
# pushing the existing extr1-extr8 rows (rows 8-15) down to rows 10-17.
$ws.Rows("8:9").Insert()

# New "line7" row (row 8)
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 1).Style = "Normal"
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New "line8" row (row 9)
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Re-number the A column (index) for the shifted extr rows (now rows 10-17) to 8..15
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15

# Update in_service flags that changed for the shifted extr rows
$ws.Cells.Item(10, 5).Value = $true   # extr1 in_service: False -> True
$ws.Cells.Item(13, 5).Value = $true   # extr4 in_service: False -> True
$ws.Cells.Item(14, 5).Value = $false  # extr5 in_service: True -> False
